# Update the "dSF" column (column F) for a set of rows with newly
# repulled data / recalculated means, per commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -2
    6  = -1
    8  = -2
    13 = 3
    15 = -1
    23 = 1
    33 = -2
    34 = 3
    38 = 4
    39 = 4
    40 = -1
    43 = 2
    50 = -3
    51 = 0
    53 = -1
    56 = 2
    60 = -3
    62 = -3
    69 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
